# The deck's two theme parts (ppt/theme/theme1.xml "Office Theme" and
# ppt/theme/theme2.xml "Integral") swap their color palettes: the theme
# actually driving the slide master/slides (reachable here through
# SlideMaster.Theme.ThemeColorScheme) goes from the "Integral" palette to
# the plain "Office Theme" palette, while the Office Theme palette moves to
# the notes-master-only theme part.
#
# fontScheme/fmtScheme are identical between the two themes already, so the
# only meaningful content change reachable through the object model is the
# 12-slot colour scheme on the presentation's (single) theme.

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

function Set-ThemeColor {
    param($Scheme, [int]$Index, [int]$R, [int]$G, [int]$B)
    $Scheme.Item($Index).RGB = $R + ($G * 256) + ($B * 65536)
}

# Target palette = the "Office Theme" colours (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink), in ThemeColorScheme's 1-12 slot order.
Set-ThemeColor $cs 1  0x00 0x00 0x00   # dk1      000000
Set-ThemeColor $cs 2  0xFF 0xFF 0xFF   # lt1      FFFFFF
Set-ThemeColor $cs 3  0x44 0x54 0x6A   # dk2      44546A
Set-ThemeColor $cs 4  0xE7 0xE6 0xE6   # lt2      E7E6E6
Set-ThemeColor $cs 5  0x5B 0x9B 0xD5   # accent1  5B9BD5
Set-ThemeColor $cs 6  0xED 0x7D 0x31   # accent2  ED7D31
Set-ThemeColor $cs 7  0xA5 0xA5 0xA5   # accent3  A5A5A5
Set-ThemeColor $cs 8  0xFF 0xC0 0x00   # accent4  FFC000
Set-ThemeColor $cs 9  0x44 0x72 0xC4   # accent5  4472C4
Set-ThemeColor $cs 10 0x70 0xAD 0x47   # accent6  70AD47
Set-ThemeColor $cs 11 0x05 0x63 0xC1   # hlink    0563C1
Set-ThemeColor $cs 12 0x95 0x4F 0x72   # folHlink 954F72
